# Auto-generated Excel COM-interop script to apply the Kujata_Profits market-data refresh.
# Updates static price/profit columns (H-N) for specific leve rows across all 8 job sheets,
# matching the values captured by the scheduled Universalis price-refresh diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 545.9091
$ws.Range("I33").Value = 552.5
$ws.Range("K33").Value = 552.5
$ws.Range("M33").Value = -323.5

$ws.Range("H74").Value = 3996.5
$ws.Range("I74").Value = 3999
$ws.Range("K74").Value = 3999
$ws.Range("M74").Value = -3063

$ws.Range("H77").Value = 3996.5
$ws.Range("I77").Value = 3999
$ws.Range("K77").Value = 19995
$ws.Range("M77").Value = -15315

$ws.Range("H112").Value = 2051.4583
$ws.Range("I112").Value = 786.55554
$ws.Range("K112").Value = 2359.66662
$ws.Range("M112").Value = -1251.66662

$ws.Range("H135").Value = 52632520
$ws.Range("I135").Value = 422.84616
$ws.Range("J135").Value = 166668740
$ws.Range("K135").Value = 3805.61544
$ws.Range("L135").Value = 1500018660
$ws.Range("M135").Value = -1270.61544
$ws.Range("N135").Value = -1500023730

$ws.Range("H138").Value = 1381.7727
$ws.Range("I138").Value = 728.23254
$ws.Range("J138").Value = 2006.2667
$ws.Range("K138").Value = 2184.69762
$ws.Range("L138").Value = 6018.800099999999
$ws.Range("M138").Value = 2955.30238
$ws.Range("N138").Value = -16298.8001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2586.6
$ws.Range("I32").Value = 2572.3572
$ws.Range("J32").Value = 2675.2222
$ws.Range("K32").Value = 2572.3572
$ws.Range("L32").Value = 2675.2222
$ws.Range("M32").Value = -2285.3572
$ws.Range("N32").Value = -3249.2222

$ws.Range("H63").Value = 1500
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 1500
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 1500
$ws.Range("N63").Value = -2872
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 1500
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 1500
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 7500
$ws.Range("N66").Value = -14364
$ws.Range("M66").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1525.2858
$ws.Range("I20").Value = 1236.5714
$ws.Range("J20").Value = 1814
$ws.Range("K20").Value = 1236.5714
$ws.Range("L20").Value = 1814
$ws.Range("M20").Value = -989.5714
$ws.Range("N20").Value = -2308

$ws.Range("H80").Value = 680.05
$ws.Range("I80").Value = 446.4
$ws.Range("J80").Value = 913.7
$ws.Range("K80").Value = 446.4
$ws.Range("L80").Value = 913.7
$ws.Range("M80").Value = 551.6
$ws.Range("N80").Value = -2909.7

$ws.Range("H83").Value = 680.05
$ws.Range("I83").Value = 446.4
$ws.Range("J83").Value = 913.7
$ws.Range("K83").Value = 2232
$ws.Range("L83").Value = 4568.5
$ws.Range("M83").Value = 2760
$ws.Range("N83").Value = -14552.5

$ws.Range("H86").Value = 3500.5938
$ws.Range("I86").Value = 3560.3704
$ws.Range("J86").Value = 3177.8
$ws.Range("K86").Value = 3560.3704
$ws.Range("L86").Value = 3177.8
$ws.Range("M86").Value = -2437.3704
$ws.Range("N86").Value = -5423.8

$ws.Range("H89").Value = 3500.5938
$ws.Range("I89").Value = 3560.3704
$ws.Range("J89").Value = 3177.8
$ws.Range("K89").Value = 17801.852
$ws.Range("L89").Value = 15889
$ws.Range("M89").Value = -12185.852
$ws.Range("N89").Value = -27121

$ws.Range("H134").Value = 4580.657
$ws.Range("I134").Value = 1159.625
$ws.Range("J134").Value = 12044.728
$ws.Range("K134").Value = 3478.875
$ws.Range("L134").Value = 36134.18399999999
$ws.Range("M134").Value = -943.875
$ws.Range("N134").Value = -41204.18399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1129.228
$ws.Range("I31").Value = 1106.6545
$ws.Range("J31").Value = 1750
$ws.Range("K31").Value = 1106.6545
$ws.Range("L31").Value = 1750
$ws.Range("M31").Value = -811.6545000000001
$ws.Range("N31").Value = -2340

$ws.Range("H34").Value = 1129.228
$ws.Range("I34").Value = 1106.6545
$ws.Range("J34").Value = 1750
$ws.Range("K34").Value = 1106.6545
$ws.Range("L34").Value = 1750
$ws.Range("M34").Value = -904.6545000000001
$ws.Range("N34").Value = -2154

$ws.Range("H62").Value = 7694595.5
$ws.Range("I62").Value = 2395
$ws.Range("K62").Value = 2395
$ws.Range("M62").Value = -1771

$ws.Range("H64").Value = 25333.334
$ws.Range("J64").Value = 25333.334
$ws.Range("L64").Value = 25333.334
$ws.Range("N64").Value = -25829.334

$ws.Range("H65").Value = 7694595.5
$ws.Range("I65").Value = 2395
$ws.Range("K65").Value = 11975
$ws.Range("M65").Value = -8855

$ws.Range("H67").Value = 25333.334
$ws.Range("J67").Value = 25333.334
$ws.Range("L67").Value = 25333.334
$ws.Range("N67").Value = -27049.334

$ws.Range("H122").Value = 872.8333
$ws.Range("I122").Value = 808.61536
$ws.Range("K122").Value = 2425.84608
$ws.Range("M122").Value = 24.15391999999974

$ws.Range("H132").Value = 3032
$ws.Range("I132").Value = 1572.5
$ws.Range("K132").Value = 4717.5
$ws.Range("M132").Value = -2187.5

$ws.Range("H134").Value = 1915.375
$ws.Range("I134").Value = 1720.5
$ws.Range("K134").Value = 5161.5
$ws.Range("M134").Value = -2626.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 250000
$ws.Range("J37").Value = 250000
$ws.Range("L37").Value = 750000
$ws.Range("N37").Value = -750224

$ws.Range("H44").Value = 1399.8667
$ws.Range("I44").Value = 1062.25
$ws.Range("J44").Value = 1785.7142
$ws.Range("K44").Value = 3186.75
$ws.Range("L44").Value = 5357.142599999999
$ws.Range("M44").Value = -2788.75
$ws.Range("N44").Value = -6153.142599999999

$ws.Range("H55").Value = 2312.25
$ws.Range("J55").Value = 3199.6
$ws.Range("L55").Value = 9598.799999999999
$ws.Range("N55").Value = -9952.799999999999

$ws.Range("H130").Value = 2018.3334
$ws.Range("J130").Value = 2018.3334
$ws.Range("L130").Value = 6055.0002
$ws.Range("N130").Value = -16095.0002

$ws.Range("H131").Value = 18521362
$ws.Range("J131").Value = 3260.913
$ws.Range("L131").Value = 9782.739
$ws.Range("N131").Value = -19862.739

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 505
$ws.Range("I97").Value = 505
$ws.Range("K97").Value = 505
$ws.Range("M97").Value = -9

$ws.Range("H132").Value = 2228.4814
$ws.Range("I132").Value = 1519.75
$ws.Range("K132").Value = 4559.25
$ws.Range("M132").Value = -2029.25

$ws.Range("H135").Value = 29666.334
$ws.Range("J135").Value = 29666.334
$ws.Range("L135").Value = 29666.334
$ws.Range("N135").Value = -39806.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1887.25
$ws.Range("I46").Value = 824.75
$ws.Range("J46").Value = 2949.75
$ws.Range("K46").Value = 824.75
$ws.Range("L46").Value = 2949.75
$ws.Range("M46").Value = -636.75
$ws.Range("N46").Value = -3325.75

$ws.Range("H61").Value = 1549.4166
$ws.Range("I61").Value = 1298.1111
$ws.Range("K61").Value = 1298.1111
$ws.Range("M61").Value = -1096.1111

$ws.Range("H93").Value = 702.6667
$ws.Range("I93").Value = 702.6667
$ws.Range("K93").Value = 702.6667
$ws.Range("M93").Value = 545.3333

$ws.Range("H113").Value = 1549.4166
$ws.Range("I113").Value = 1298.1111
$ws.Range("K113").Value = 1298.1111
$ws.Range("M113").Value = 871.8888999999999

$ws.Range("H132").Value = 29454.695
$ws.Range("I132").Value = 1070.9546
$ws.Range("J132").Value = 74057.71000000001
$ws.Range("K132").Value = 3212.8638
$ws.Range("L132").Value = 222173.13
$ws.Range("M132").Value = -682.8638000000001
$ws.Range("N132").Value = -227233.13

$ws.Range("H136").Value = 1416.5883
$ws.Range("I136").Value = 1385.1333
$ws.Range("J136").Value = 1652.5
$ws.Range("K136").Value = 4155.3999
$ws.Range("L136").Value = 4957.5
$ws.Range("M136").Value = -1605.3999
$ws.Range("N136").Value = -10057.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 10124.5
$ws.Range("J63").Value = 10124.5
$ws.Range("L63").Value = 10124.5
$ws.Range("N63").Value = -11372.5

$ws.Range("H66").Value = 10124.5
$ws.Range("J66").Value = 10124.5
$ws.Range("L66").Value = 30373.5
$ws.Range("N66").Value = -36613.5

$ws.Range("H136").Value = 811.35297
$ws.Range("I136").Value = 579.8333
$ws.Range("K136").Value = 1739.4999
$ws.Range("M136").Value = 810.5001
